# Update evidence for task B1-B4
$wb = $excel.ActiveWorkbook

# --- B1: record the first pair of new hash values ---
$wsB1 = $wb.Worksheets.Item("B1")
$wsB1.Range("A2").Value = "E8750EF477359BFC9C9E6271E1DE16ECFF0D1AC4DE27B5D649C586AF0FCB32D7"
$wsB1.Range("A3").Value = "03E4D08C44032891A8577581B10D5899CADD768B4C8D6BD1606AD76CE9032005"
$null = $wsB1.Range("A3").Select()

# --- B2: record the second pair of new hash values ---
$wsB2 = $wb.Worksheets.Item("B2")
$wsB2.Range("A2").Value = "6FC5A35A43506803589761F04554EDF5009C4F4FED438FB782CB7899EC23F50C"
$wsB2.Range("A3").Value = "643454247835F2A7F7CE4F251E55C009B150A220836D8A733A7AE96B47FFE416"

# B2 becomes the active sheet/tab, with A3 selected
$null = $wsB2.Activate()
$null = $wsB2.Range("A3").Select()
